$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-affecting data rows 8-37 with new values from the diff.
# Row 8
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = ""
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 38
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "1.0"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it's ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "0.00"
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).NumberFormat = "@"
$ws.Cells.Item(8, 9).Value = ""
# Row 9
$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = "P. point"
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = 70
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "2"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "Short point (up to 3 mtr.)"
$ws.Cells.Item(9, 6).Value = 256
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "17920.00"
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).NumberFormat = "@"
$ws.Cells.Item(9, 9).Value = ""
# Row 10
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "P. point"
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 41
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "3"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "Medium point (up to 6 mtr.)"
$ws.Cells.Item(10, 6).Value = 472
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = "19352.00"
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).NumberFormat = "@"
$ws.Cells.Item(10, 9).Value = ""
# Row 11
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "P. point"
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 85
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "4"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "Long point  (up to 10 mtr.)"
$ws.Cells.Item(11, 6).Value = 662
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = "56270.00"
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).NumberFormat = "@"
$ws.Cells.Item(11, 9).Value = ""
# Row 12
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "P. point"
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 47
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "6"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "On board"
$ws.Cells.Item(12, 6).Value = 136
$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = "6392.00"
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).NumberFormat = "@"
$ws.Cells.Item(12, 9).Value = ""
# Row 13
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "Each"
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 70
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "3.0"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Cells.Item(13, 6).Value = 23
$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = "1610.00"
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).NumberFormat = "@"
$ws.Cells.Item(13, 9).Value = ""
# Row 14
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Value = "Each"
$ws.Cells.Item(14, 2).Value = 0
$ws.Cells.Item(14, 3).Value = 23
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.0"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Cells.Item(14, 6).Value = 50
$ws.Cells.Item(14, 7).NumberFormat = "@"
$ws.Cells.Item(14, 7).Value = "1150.00"
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).NumberFormat = "@"
$ws.Cells.Item(14, 9).Value = ""
# Row 15
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "Each"
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(15, 3).Value = 22
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "5.0"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "Providing & Fixing of  of 3/5 pin 6 amp. flush type  non modular socket  made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Cells.Item(15, 6).Value = 33
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "726.00"
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).NumberFormat = "@"
$ws.Cells.Item(15, 9).Value = ""
# Row 16
$ws.Cells.Item(16, 1).NumberFormat = "@"
$ws.Cells.Item(16, 1).Value = "Each"
$ws.Cells.Item(16, 2).Value = 0
$ws.Cells.Item(16, 3).Value = 70
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "6.0"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Cells.Item(16, 6).Value = 78
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "5460.00"
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).NumberFormat = "@"
$ws.Cells.Item(16, 9).Value = ""
# Row 17
$ws.Cells.Item(17, 1).NumberFormat = "@"
$ws.Cells.Item(17, 1).Value = "Each"
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(17, 3).Value = 21
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "9.0"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Cells.Item(17, 6).Value = 219
$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "4599.00"
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).NumberFormat = "@"
$ws.Cells.Item(17, 9).Value = ""
# Row 18
$ws.Cells.Item(18, 1).NumberFormat = "@"
$ws.Cells.Item(18, 1).Value = "Each"
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(18, 3).Value = 35
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "10.0"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Cells.Item(18, 6).Value = 303
$ws.Cells.Item(18, 7).NumberFormat = "@"
$ws.Cells.Item(18, 7).Value = "10605.00"
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).NumberFormat = "@"
$ws.Cells.Item(18, 9).Value = ""
# Row 19
$ws.Cells.Item(19, 1).NumberFormat = "@"
$ws.Cells.Item(19, 1).Value = ""
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 53
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "11.0"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).NumberFormat = "@"
$ws.Cells.Item(19, 7).Value = "0.00"
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).NumberFormat = "@"
$ws.Cells.Item(19, 9).Value = ""
# Row 20
$ws.Cells.Item(20, 1).NumberFormat = "@"
$ws.Cells.Item(20, 1).Value = "R. mtr."
$ws.Cells.Item(20, 2).Value = 0
$ws.Cells.Item(20, 3).Value = 39
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "16"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "20 mm"
$ws.Cells.Item(20, 6).Value = 40
$ws.Cells.Item(20, 7).NumberFormat = "@"
$ws.Cells.Item(20, 7).Value = "1560.00"
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).NumberFormat = "@"
$ws.Cells.Item(20, 9).Value = ""
# Row 21
$ws.Cells.Item(21, 1).NumberFormat = "@"
$ws.Cells.Item(21, 1).Value = "R. mtr."
$ws.Cells.Item(21, 2).Value = 0
$ws.Cells.Item(21, 3).Value = 55
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "17"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "25 mm"
$ws.Cells.Item(21, 6).Value = 56
$ws.Cells.Item(21, 7).NumberFormat = "@"
$ws.Cells.Item(21, 7).Value = "3080.00"
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).NumberFormat = "@"
$ws.Cells.Item(21, 9).Value = ""
# Row 22
$ws.Cells.Item(22, 1).NumberFormat = "@"
$ws.Cells.Item(22, 1).Value = ""
$ws.Cells.Item(22, 2).Value = 0
$ws.Cells.Item(22, 3).Value = 20
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "14.0"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR .   "
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).NumberFormat = "@"
$ws.Cells.Item(22, 7).Value = "0.00"
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).NumberFormat = "@"
$ws.Cells.Item(22, 9).Value = ""
# Row 23
$ws.Cells.Item(23, 1).NumberFormat = "@"
$ws.Cells.Item(23, 1).Value = "Mtr."
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 52
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "23"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "8 SWG G.I. ( Hot Dipped  ) Wire "
$ws.Cells.Item(23, 6).Value = 20
$ws.Cells.Item(23, 7).NumberFormat = "@"
$ws.Cells.Item(23, 7).Value = "1040.00"
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).NumberFormat = "@"
$ws.Cells.Item(23, 9).Value = ""
# Row 24
$ws.Cells.Item(24, 1).NumberFormat = "@"
$ws.Cells.Item(24, 1).Value = ""
$ws.Cells.Item(24, 2).Value = 0
$ws.Cells.Item(24, 3).Value = 3
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "15.0"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).NumberFormat = "@"
$ws.Cells.Item(24, 7).Value = "0.00"
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 9).NumberFormat = "@"
$ws.Cells.Item(24, 9).Value = ""
# Row 25
$ws.Cells.Item(25, 1).NumberFormat = "@"
$ws.Cells.Item(25, 1).Value = "Each"
$ws.Cells.Item(25, 2).Value = 0
$ws.Cells.Item(25, 3).Value = 54
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "25"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )"
$ws.Cells.Item(25, 6).Value = 1890
$ws.Cells.Item(25, 7).NumberFormat = "@"
$ws.Cells.Item(25, 7).Value = "102060.00"
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).NumberFormat = "@"
$ws.Cells.Item(25, 9).Value = ""
# Row 26
$ws.Cells.Item(26, 1).NumberFormat = "@"
$ws.Cells.Item(26, 1).Value = ""
$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(26, 3).Value = 99
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "29"
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "Single pole MCB   (With B/C curve tripping Characteristics)"
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).NumberFormat = "@"
$ws.Cells.Item(26, 7).Value = "0.00"
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).NumberFormat = "@"
$ws.Cells.Item(26, 9).Value = ""
# Row 27
$ws.Cells.Item(27, 1).NumberFormat = "@"
$ws.Cells.Item(27, 1).Value = "Each"
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(27, 3).Value = 37
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "32"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = " 50/63 A rating"
$ws.Cells.Item(27, 6).Value = 900
$ws.Cells.Item(27, 7).NumberFormat = "@"
$ws.Cells.Item(27, 7).Value = "33300.00"
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).NumberFormat = "@"
$ws.Cells.Item(27, 9).Value = ""
# Row 28
$ws.Cells.Item(28, 1).NumberFormat = "@"
$ws.Cells.Item(28, 1).Value = ""
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(28, 3).Value = 15
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "18.0"
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).NumberFormat = "@"
$ws.Cells.Item(28, 7).Value = "0.00"
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).NumberFormat = "@"
$ws.Cells.Item(28, 9).Value = ""
# Row 29
$ws.Cells.Item(29, 1).NumberFormat = "@"
$ws.Cells.Item(29, 1).Value = ""
$ws.Cells.Item(29, 2).Value = 0
$ws.Cells.Item(29, 3).Value = 69
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "34"
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = "Metal door (single phase) IK-09 and IP-43 with Metal end box"
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).NumberFormat = "@"
$ws.Cells.Item(29, 7).Value = "0.00"
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).NumberFormat = "@"
$ws.Cells.Item(29, 9).Value = ""
# Row 30
$ws.Cells.Item(30, 1).NumberFormat = "@"
$ws.Cells.Item(30, 1).Value = "Each"
$ws.Cells.Item(30, 2).Value = 0
$ws.Cells.Item(30, 3).Value = 77
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "35"
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = "8 Way (8+2)"
$ws.Cells.Item(30, 6).Value = 2184
$ws.Cells.Item(30, 7).NumberFormat = "@"
$ws.Cells.Item(30, 7).Value = "168168.00"
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).NumberFormat = "@"
$ws.Cells.Item(30, 9).Value = ""
# Row 31
$ws.Cells.Item(31, 1).NumberFormat = "@"
$ws.Cells.Item(31, 1).Value = ""
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 91
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "36"
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = "Total"
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).NumberFormat = "@"
$ws.Cells.Item(31, 7).Value = "0.00"
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).NumberFormat = "@"
$ws.Cells.Item(31, 9).Value = ""
# Row 32
$ws.Cells.Item(32, 1).NumberFormat = "@"
$ws.Cells.Item(32, 1).Value = "%"
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(32, 3).Value = 31
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "37"
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = "Add Tender Premium "
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).NumberFormat = "@"
$ws.Cells.Item(32, 7).Value = "0.00"
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).NumberFormat = "@"
$ws.Cells.Item(32, 9).Value = ""
# Row 33
$ws.Cells.Item(33, 1).NumberFormat = "@"
$ws.Cells.Item(33, 1).Value = ""
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 36
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "38"
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = "Grand Total"
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).NumberFormat = "@"
$ws.Cells.Item(33, 7).Value = "0.00"
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).NumberFormat = "@"
$ws.Cells.Item(33, 9).Value = ""
# Row 34
$ws.Cells.Item(34, 1).NumberFormat = "@"
$ws.Cells.Item(34, 1).Value = ""
# Row 35
$ws.Cells.Item(35, 1).NumberFormat = "@"
$ws.Cells.Item(35, 1).Value = ""
$ws.Cells.Item(35, 2).NumberFormat = "@"
$ws.Cells.Item(35, 2).Value = ""
$ws.Cells.Item(35, 3).NumberFormat = "@"
$ws.Cells.Item(35, 3).Value = ""
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = ""
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = "Grand Total Rs."
$ws.Cells.Item(35, 6).NumberFormat = "@"
$ws.Cells.Item(35, 6).Value = ""
$ws.Cells.Item(35, 7).NumberFormat = "@"
$ws.Cells.Item(35, 7).Value = "433292.00"
$ws.Cells.Item(35, 8).NumberFormat = "@"
$ws.Cells.Item(35, 8).Value = "433292.00"
$ws.Cells.Item(35, 9).NumberFormat = "@"
$ws.Cells.Item(35, 9).Value = ""
# Row 36
$ws.Cells.Item(36, 1).NumberFormat = "@"
$ws.Cells.Item(36, 1).Value = ""
$ws.Cells.Item(36, 2).NumberFormat = "@"
$ws.Cells.Item(36, 2).Value = ""
$ws.Cells.Item(36, 3).NumberFormat = "@"
$ws.Cells.Item(36, 3).Value = ""
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = ""
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = "Tender Premium @ 0%"
$ws.Cells.Item(36, 6).NumberFormat = "@"
$ws.Cells.Item(36, 6).Value = ""
$ws.Cells.Item(36, 7).NumberFormat = "@"
$ws.Cells.Item(36, 7).Value = "0.00"
$ws.Cells.Item(36, 8).NumberFormat = "@"
$ws.Cells.Item(36, 8).Value = "0.00"
$ws.Cells.Item(36, 9).NumberFormat = "@"
$ws.Cells.Item(36, 9).Value = ""
# Row 37
$ws.Cells.Item(37, 1).NumberFormat = "@"
$ws.Cells.Item(37, 1).Value = ""
$ws.Cells.Item(37, 2).NumberFormat = "@"
$ws.Cells.Item(37, 2).Value = ""
$ws.Cells.Item(37, 3).NumberFormat = "@"
$ws.Cells.Item(37, 3).Value = ""
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = ""
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = "NET PAYABLE AMOUNT Rs."
$ws.Cells.Item(37, 6).NumberFormat = "@"
$ws.Cells.Item(37, 6).Value = ""
$ws.Cells.Item(37, 7).NumberFormat = "@"
$ws.Cells.Item(37, 7).Value = "433292.00"
$ws.Cells.Item(37, 8).NumberFormat = "@"
$ws.Cells.Item(37, 8).Value = "433292.00"
$ws.Cells.Item(37, 9).NumberFormat = "@"
$ws.Cells.Item(37, 9).Value = ""

# Row 34 loses its data in columns B:I (only column A, already empty, remains).
$ws.Range("B34:I34").ClearContents()

# Row 38 is entirely removed from the used range (the table now ends at row 37).
$ws.Range("A38:I38").ClearContents()

# Shrink the sheet dimension to match the new extent.
$ws.UsedRange | Out-Null